# Apply odds updates per diff for Jogos_da_Semana_FlashScore_2024-10-12.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("J3").Value = 2.38
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 6.5
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("S3").Value = 1.53
$ws.Range("T3").Value = 2.38
$ws.Range("W3").Value = 5
$ws.Range("AC3").Value = 6.5
$ws.Range("AL3").Value = 67
$ws.Range("AP3").Value = 26
$ws.Range("AQ3").Value = 34
$ws.Range("AR3").Value = 67
$ws.Range("AT3").Value = 2.38

# Row 11
$ws.Range("N11").Value = 6.5

# Row 26
$ws.Range("Z26").Value = 11
$ws.Range("AM26").Value = 126
$ws.Range("AR26").Value = 29
$ws.Range("AX26").Value = 29

# Row 27
$ws.Range("G27").Value = 2.22
$ws.Range("H27").Value = 3.15
$ws.Range("I27").Value = 3

# Row 28
$ws.Range("G28").Value = 2.75
$ws.Range("H28").Value = 3.15
$ws.Range("I28").Value = 2.37

# Row 29
$ws.Range("G29").Value = 2.5
$ws.Range("H29").Value = 3.4
$ws.Range("I29").Value = 2.45

# Row 30
$ws.Range("H30").Value = 4.45
$ws.Range("I30").Value = 4.9

# Row 31
$ws.Range("H31").Value = 3.6
$ws.Range("I31").Value = 4

# Row 32
$ws.Range("G32").Value = 2.55
$ws.Range("H32").Value = 3.15
$ws.Range("I32").Value = 2.55

# Row 33
$ws.Range("G33").Value = 3.2
$ws.Range("H33").Value = 3.6
$ws.Range("I33").Value = 1.98

# Row 34
$ws.Range("G34").Value = 2.55
$ws.Range("H34").Value = 3.35
$ws.Range("I34").Value = 2.45

# Row 35
$ws.Range("G35").Value = 1.78
$ws.Range("H35").Value = 3.4
$ws.Range("I35").Value = 4.3
$ws.Range("J35").Value = 2.32
$ws.Range("K35").Value = 2.12
$ws.Range("L35").Value = 4.5
$ws.Range("M35").Value = 1.01
$ws.Range("N35").Value = 8.1
$ws.Range("O35").Value = 1.3
$ws.Range("P35").Value = 2.92
$ws.Range("Q35").Value = 1.93
$ws.Range("R35").Value = 1.78
$ws.Range("U35").Value = 1.8
$ws.Range("V35").Value = 1.8
$ws.Range("W35").Value = 6.5
$ws.Range("X35").Value = 8
$ws.Range("Y35").Value = 8.25
$ws.Range("Z35").Value = 14.5
$ws.Range("AA35").Value = 15
$ws.Range("AB35").Value = 28
$ws.Range("AC35").Value = 9.25
$ws.Range("AD35").Value = 6.6
$ws.Range("AE35").Value = 15.5
$ws.Range("AF35").Value = 75
$ws.Range("AG35").Value = 11.75
$ws.Range("AH35").Value = 24
$ws.Range("AI35").Value = 14
$ws.Range("AJ35").Value = 70
$ws.Range("AK35").Value = 40
$ws.Range("AM35").Value = 600
$ws.Range("AN35").Value = 3.6
$ws.Range("AO35").Value = 8.75
$ws.Range("AP35").Value = 17.5
$ws.Range("AQ35").Value = 30
$ws.Range("AR35").Value = 60
$ws.Range("AS35").Value = 250
$ws.Range("AT35").Value = 2.6
$ws.Range("AU35").Value = 7.1
$ws.Range("AV35").Value = 60
$ws.Range("AW35").Value = 6
$ws.Range("AX35").Value = 23
$ws.Range("AY35").Value = 28
$ws.Range("AZ35").Value = 120
$ws.Range("BA35").Value = 150
$ws.Range("BB35").Value = 350
